$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update time_taken (column F) values in the "data" sheet
$ws.Cells.Item(2, 6).Value = "2021-10-05 14:19:59.545495"
$ws.Cells.Item(3, 6).Value = "2021-10-05 14:19:59.545503"
$ws.Cells.Item(4, 6).Value = "2021-10-05 14:19:59.545507"
$ws.Cells.Item(5, 6).Value = "2021-10-05 14:19:59.545510"
$ws.Cells.Item(6, 6).Value = "2021-10-05 14:19:59.545513"
$ws.Cells.Item(7, 6).Value = "2021-10-05 14:19:59.545515"
$ws.Cells.Item(8, 6).Value = "2021-10-05 14:19:59.545518"
$ws.Cells.Item(9, 6).Value = "2021-10-05 14:19:59.545520"
$ws.Cells.Item(10, 6).Value = "2021-10-05 14:19:59.545523"
$ws.Cells.Item(11, 6).Value = "2021-10-05 14:19:59.545526"
$ws.Cells.Item(12, 6).Value = "2021-10-05 14:19:59.545528"
$ws.Cells.Item(13, 6).Value = "2021-10-05 14:19:59.545531"
$ws.Cells.Item(14, 6).Value = "2021-10-05 14:19:59.545533"
$ws.Cells.Item(15, 6).Value = "2021-10-05 14:19:59.545536"
$ws.Cells.Item(16, 6).Value = "2021-10-05 14:19:59.545538"
$ws.Cells.Item(17, 6).Value = "2021-10-05 14:19:59.545541"
$ws.Cells.Item(18, 6).Value = "2021-10-05 14:19:59.545543"
$ws.Cells.Item(19, 6).Value = "2021-10-05 14:19:59.545546"
$ws.Cells.Item(20, 6).Value = "2021-10-05 14:19:59.545549"
$ws.Cells.Item(21, 6).Value = "2021-10-05 14:19:59.545552"
$ws.Cells.Item(22, 6).Value = "2021-10-05 14:19:59.545554"
$ws.Cells.Item(23, 6).Value = "2021-10-05 14:19:59.545557"
$ws.Cells.Item(24, 6).Value = "2021-10-05 14:19:59.545559"
$ws.Cells.Item(25, 6).Value = "2021-10-05 14:19:59.545562"
$ws.Cells.Item(26, 6).Value = "2021-10-05 14:19:59.545564"
$ws.Cells.Item(27, 6).Value = "2021-10-05 14:19:59.545567"
$ws.Cells.Item(28, 6).Value = "2021-10-05 14:19:59.545570"
$ws.Cells.Item(29, 6).Value = "2021-10-05 14:19:59.545573"
$ws.Cells.Item(30, 6).Value = "2021-10-05 14:19:59.545575"
$ws.Cells.Item(31, 6).Value = "2021-10-05 14:19:59.545578"
$ws.Cells.Item(32, 6).Value = "2021-10-05 14:19:59.545580"
$ws.Cells.Item(33, 6).Value = "2021-10-05 14:19:59.545582"
$ws.Cells.Item(34, 6).Value = "2021-10-05 14:19:59.545585"
$ws.Cells.Item(35, 6).Value = "2021-10-05 14:19:59.545591"
$ws.Cells.Item(36, 6).Value = "2021-10-05 14:19:59.545594"
$ws.Cells.Item(37, 6).Value = "2021-10-05 14:19:59.545596"
$ws.Cells.Item(38, 6).Value = "2021-10-05 14:19:59.545599"
$ws.Cells.Item(39, 6).Value = "2021-10-05 14:19:59.545601"
$ws.Cells.Item(40, 6).Value = "2021-10-05 14:19:59.545604"
$ws.Cells.Item(41, 6).Value = "2021-10-05 14:19:59.545606"
$ws.Cells.Item(42, 6).Value = "2021-10-05 14:19:59.545609"
$ws.Cells.Item(43, 6).Value = "2021-10-05 14:19:59.545612"
$ws.Cells.Item(44, 6).Value = "2021-10-05 14:19:59.545615"
$ws.Cells.Item(45, 6).Value = "2021-10-05 14:19:59.545617"
$ws.Cells.Item(46, 6).Value = "2021-10-05 14:19:59.545620"
$ws.Cells.Item(47, 6).Value = "2021-10-05 14:19:59.545622"
$ws.Cells.Item(48, 6).Value = "2021-10-05 14:19:59.545625"
$ws.Cells.Item(49, 6).Value = "2021-10-05 14:19:59.545628"
$ws.Cells.Item(50, 6).Value = "2021-10-05 14:19:59.545630"
$ws.Cells.Item(51, 6).Value = "2021-10-05 14:19:59.545633"
$ws.Cells.Item(52, 6).Value = "2021-10-05 14:19:59.545635"
$ws.Cells.Item(53, 6).Value = "2021-10-05 14:19:59.545638"
$ws.Cells.Item(54, 6).Value = "2021-10-05 14:19:59.545641"
$ws.Cells.Item(55, 6).Value = "2021-10-05 14:19:59.545644"
$ws.Cells.Item(56, 6).Value = "2021-10-05 14:19:59.545646"
$ws.Cells.Item(57, 6).Value = "2021-10-05 14:19:59.545649"
$ws.Cells.Item(58, 6).Value = "2021-10-05 14:19:59.545651"
$ws.Cells.Item(59, 6).Value = "2021-10-05 14:19:59.545654"
$ws.Cells.Item(60, 6).Value = "2021-10-05 14:19:59.545656"
$ws.Cells.Item(61, 6).Value = "2021-10-05 14:19:59.545659"
$ws.Cells.Item(62, 6).Value = "2021-10-05 14:19:59.545662"
$ws.Cells.Item(63, 6).Value = "2021-10-05 14:19:59.545664"
$ws.Cells.Item(64, 6).Value = "2021-10-05 14:19:59.545667"
$ws.Cells.Item(65, 6).Value = "2021-10-05 14:19:59.545670"
# Add the new "metadata" worksheet as the last sheet
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$metaSheet.Name = "metadata"

# Header row (row 1), bold/boxed header style matches the "data" sheet's header style
$metaSheet.Cells.Item(1, 2).Value = "data_name"
$metaSheet.Cells.Item(1, 3).Value = "data_id"
$metaSheet.Cells.Item(1, 4).Value = "data_version"
$metaSheet.Cells.Item(1, 5).Value = "data_version_created"
$metaSheet.Cells.Item(1, 6).Value = "panel_query_time"
$metaSheet.Cells.Item(1, 7).Value = "panel_get_request"

# Data row (row 2)
$metaSheet.Cells.Item(2, 1).Value = 0
$metaSheet.Cells.Item(2, 2).Value = "Diabetes with additional phenotypes suggestive of a monogenic aetiology"
$metaSheet.Cells.Item(2, 3).Value = 26

# data_version must remain a text value ("1.63"), not be auto-converted to a number
$metaSheet.Cells.Item(2, 4).NumberFormat = "@"
$metaSheet.Cells.Item(2, 4).Value = "1.63"
$metaSheet.Cells.Item(2, 4).Style = "Normal"

$metaSheet.Cells.Item(2, 5).Value = "2021-07-28T09:54:50.314484Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:19:59.542133"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/26/?format=json"

# Apply the same header formatting used on the "data" sheet's header row
# (bold font, thin box border, centered horizontal / top vertical alignment)
$headerRange = $metaSheet.Range($metaSheet.Cells.Item(1, 2), $metaSheet.Cells.Item(1, 7))
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$indexCell = $metaSheet.Cells.Item(2, 1)
$indexCell.Font.Bold = $true
$indexCell.HorizontalAlignment = -4108
$indexCell.VerticalAlignment = -4160
$indexCell.Borders.LineStyle = 1

# Make sure the "data" sheet remains the active tab, as in the original workbook
$ws.Activate()
